$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Value2

    if ($text -ne $null -and $text -ne "") {
        $parts = $text -split ", "
        if ($parts.Length -gt 1 -and $parts[0] -eq "System") {
            $n = $parts.Length
            $reversed = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $newText = $reversed -join ", "
            $cell.Value = $newText
        }
    }
}
